# Apply edits for 06_pin_tot_by_admin_area_severity.xlsx
# Source data re-aggregated from admin3 -> admin1 level (fewer, broader rows).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell A1: admin3 -> admin1
$ws.Cells.Item(1,1).Value = "admin1"

# Row 2: MMR001
$ws.Cells.Item(2,1).Value = "MMR001"
$ws.Cells.Item(2,2).Value = 114297
$ws.Cells.Item(2,3).Value = 72.5
$ws.Cells.Item(2,4).Value = 82816
$ws.Cells.Item(2,5).Value = 23.8
$ws.Cells.Item(2,6).Value = 27148
$ws.Cells.Item(2,7).Value = 3.2
$ws.Cells.Item(2,8).Value = 3711
$ws.Cells.Item(2,9).Value = 0.5
$ws.Cells.Item(2,10).Value = 621
$ws.Cells.Item(2,11).Value = 27.5
$ws.Cells.Item(2,12).Value = 31481
$ws.Cells.Item(2,13).Value = "'3"

# Row 3: MMR002
$ws.Cells.Item(3,1).Value = "MMR002"
$ws.Cells.Item(3,2).Value = 59845
$ws.Cells.Item(3,3).Value = 28.6
$ws.Cells.Item(3,4).Value = 17095
$ws.Cells.Item(3,5).Value = 66.90000000000001
$ws.Cells.Item(3,6).Value = 40017
$ws.Cells.Item(3,7).Value = 3.5
$ws.Cells.Item(3,8).Value = 2112
$ws.Cells.Item(3,9).Value = 1
$ws.Cells.Item(3,10).Value = 622
$ws.Cells.Item(3,11).Value = 71.40000000000001
$ws.Cells.Item(3,12).Value = 42751
$ws.Cells.Item(3,13).Value = "'3"

# Row 4: MMR003
$ws.Cells.Item(4,1).Value = "MMR003"
$ws.Cells.Item(4,2).Value = 86336
$ws.Cells.Item(4,3).Value = 52.9
$ws.Cells.Item(4,4).Value = 45681
$ws.Cells.Item(4,5).Value = 41.2
$ws.Cells.Item(4,6).Value = 35594
$ws.Cells.Item(4,7).Value = 5
$ws.Cells.Item(4,8).Value = 4341
$ws.Cells.Item(4,9).Value = 0.8
$ws.Cells.Item(4,10).Value = 720
$ws.Cells.Item(4,11).Value = 47.1
$ws.Cells.Item(4,12).Value = 40655
$ws.Cells.Item(4,13).Value = "'3"

# Row 5: MMR004
$ws.Cells.Item(5,1).Value = "MMR004"
$ws.Cells.Item(5,2).Value = 96362
$ws.Cells.Item(5,3).Value = 29.6
$ws.Cells.Item(5,4).Value = 28517
$ws.Cells.Item(5,5).Value = 66
$ws.Cells.Item(5,6).Value = 63559
$ws.Cells.Item(5,7).Value = 3.5
$ws.Cells.Item(5,8).Value = 3334
$ws.Cells.Item(5,9).Value = 1
$ws.Cells.Item(5,10).Value = 952
$ws.Cells.Item(5,11).Value = 70.40000000000001
$ws.Cells.Item(5,12).Value = 67845
$ws.Cells.Item(5,13).Value = "'3"

# Row 6: MMR005
$ws.Cells.Item(6,1).Value = "MMR005"
$ws.Cells.Item(6,2).Value = 251037
$ws.Cells.Item(6,3).Value = 18.7
$ws.Cells.Item(6,4).Value = 46878
$ws.Cells.Item(6,5).Value = 75.7
$ws.Cells.Item(6,6).Value = 190045
$ws.Cells.Item(6,7).Value = 4.2
$ws.Cells.Item(6,8).Value = 10513
$ws.Cells.Item(6,9).Value = 1.4
$ws.Cells.Item(6,10).Value = 3601
$ws.Cells.Item(6,11).Value = 81.3
$ws.Cells.Item(6,12).Value = 204159
$ws.Cells.Item(6,13).Value = "'3"

# Row 7: MMR006
$ws.Cells.Item(7,1).Value = "MMR006"
$ws.Cells.Item(7,2).Value = 345131
$ws.Cells.Item(7,3).Value = 67
$ws.Cells.Item(7,4).Value = 231393
$ws.Cells.Item(7,5).Value = 28.5
$ws.Cells.Item(7,6).Value = 98259
$ws.Cells.Item(7,7).Value = 4.1
$ws.Cells.Item(7,8).Value = 14278
$ws.Cells.Item(7,9).Value = 0.3
$ws.Cells.Item(7,10).Value = 1201
$ws.Cells.Item(7,11).Value = 33
$ws.Cells.Item(7,12).Value = 113737
$ws.Cells.Item(7,13).Value = "'3"

# Row 8: MMR007
$ws.Cells.Item(8,1).Value = "MMR007"
$ws.Cells.Item(8,2).Value = 138589
$ws.Cells.Item(8,3).Value = 74.5
$ws.Cells.Item(8,4).Value = 103188
$ws.Cells.Item(8,5).Value = 23.2
$ws.Cells.Item(8,6).Value = 32126
$ws.Cells.Item(8,7).Value = 2.4
$ws.Cells.Item(8,8).Value = 3275
$ws.Cells.Item(8,9).Value = 0
$ws.Cells.Item(8,10).Value = 0
$ws.Cells.Item(8,11).Value = 25.5
$ws.Cells.Item(8,12).Value = 35401
$ws.Cells.Item(8,13).Value = "'3"

# Row 9: MMR008
$ws.Cells.Item(9,1).Value = "MMR008"
$ws.Cells.Item(9,2).Value = 82380
$ws.Cells.Item(9,3).Value = 75.7
$ws.Cells.Item(9,4).Value = 62340
$ws.Cells.Item(9,5).Value = 16.8
$ws.Cells.Item(9,6).Value = 13805
$ws.Cells.Item(9,7).Value = 7.6
$ws.Cells.Item(9,8).Value = 6236
$ws.Cells.Item(9,9).Value = 0
$ws.Cells.Item(9,10).Value = 0
$ws.Cells.Item(9,11).Value = 24.3
$ws.Cells.Item(9,12).Value = 20041
$ws.Cells.Item(9,13).Value = "'3"

# Row 10: MMR009
$ws.Cells.Item(10,1).Value = "MMR009"
$ws.Cells.Item(10,2).Value = 94677
$ws.Cells.Item(10,3).Value = 71.09999999999999
$ws.Cells.Item(10,4).Value = 67289
$ws.Cells.Item(10,5).Value = 25.2
$ws.Cells.Item(10,6).Value = 23845
$ws.Cells.Item(10,7).Value = 2.8
$ws.Cells.Item(10,8).Value = 2647
$ws.Cells.Item(10,9).Value = 0.9
$ws.Cells.Item(10,10).Value = 896
$ws.Cells.Item(10,11).Value = 28.9
$ws.Cells.Item(10,12).Value = 27388
$ws.Cells.Item(10,13).Value = "'3"

# Row 11: MMR010
$ws.Cells.Item(11,1).Value = "MMR010"
$ws.Cells.Item(11,2).Value = 469612
$ws.Cells.Item(11,3).Value = 77.40000000000001
$ws.Cells.Item(11,4).Value = 363402
$ws.Cells.Item(11,5).Value = 18.4
$ws.Cells.Item(11,6).Value = 86190
$ws.Cells.Item(11,7).Value = 3.5
$ws.Cells.Item(11,8).Value = 16642
$ws.Cells.Item(11,9).Value = 0.7
$ws.Cells.Item(11,10).Value = 3379
$ws.Cells.Item(11,11).Value = 22.6
$ws.Cells.Item(11,12).Value = 106210
$ws.Cells.Item(11,13).Value = "'3"

# Row 12: MMR011
$ws.Cells.Item(12,1).Value = "MMR011"
$ws.Cells.Item(12,2).Value = 113855
$ws.Cells.Item(12,3).Value = 64.7
$ws.Cells.Item(12,4).Value = 73705
$ws.Cells.Item(12,5).Value = 18.5
$ws.Cells.Item(12,6).Value = 21096
$ws.Cells.Item(12,7).Value = 8.800000000000001
$ws.Cells.Item(12,8).Value = 10036
$ws.Cells.Item(12,9).Value = 7.9
$ws.Cells.Item(12,10).Value = 9017
$ws.Cells.Item(12,11).Value = 35.3
$ws.Cells.Item(12,12).Value = 40150
$ws.Cells.Item(12,13).Value = "'3"

# Row 13: MMR012
$ws.Cells.Item(13,1).Value = "MMR012"
$ws.Cells.Item(13,2).Value = 147696
$ws.Cells.Item(13,3).Value = 41.8
$ws.Cells.Item(13,4).Value = 61772
$ws.Cells.Item(13,5).Value = 51.1
$ws.Cells.Item(13,6).Value = 75470
$ws.Cells.Item(13,7).Value = 6.6
$ws.Cells.Item(13,8).Value = 9675
$ws.Cells.Item(13,9).Value = 0.5
$ws.Cells.Item(13,10).Value = 779
$ws.Cells.Item(13,11).Value = 58.2
$ws.Cells.Item(13,12).Value = 85924
$ws.Cells.Item(13,13).Value = "'3"

# Row 14: MMR013
$ws.Cells.Item(14,1).Value = "MMR013"
$ws.Cells.Item(14,2).Value = 94050
$ws.Cells.Item(14,3).Value = 65.8
$ws.Cells.Item(14,4).Value = 61865
$ws.Cells.Item(14,5).Value = 26
$ws.Cells.Item(14,6).Value = 24464
$ws.Cells.Item(14,7).Value = 7.8
$ws.Cells.Item(14,8).Value = 7317
$ws.Cells.Item(14,9).Value = 0.4
$ws.Cells.Item(14,10).Value = 403
$ws.Cells.Item(14,11).Value = 34.2
$ws.Cells.Item(14,12).Value = 32185
$ws.Cells.Item(14,13).Value = "'3"

# Row 15: MMR014
$ws.Cells.Item(15,1).Value = "MMR014"
$ws.Cells.Item(15,2).Value = 291446
$ws.Cells.Item(15,3).Value = 77.09999999999999
$ws.Cells.Item(15,4).Value = 224720
$ws.Cells.Item(15,5).Value = 17.6
$ws.Cells.Item(15,6).Value = 51440
$ws.Cells.Item(15,7).Value = 5.2
$ws.Cells.Item(15,8).Value = 15286
$ws.Cells.Item(15,9).Value = 0
$ws.Cells.Item(15,10).Value = 0
$ws.Cells.Item(15,11).Value = 22.9
$ws.Cells.Item(15,12).Value = 66726
$ws.Cells.Item(15,13).Value = "'3"

# Row 16: MMR015
$ws.Cells.Item(16,1).Value = "MMR015"
$ws.Cells.Item(16,2).Value = 210768
$ws.Cells.Item(16,3).Value = 48.4
$ws.Cells.Item(16,4).Value = 101932
$ws.Cells.Item(16,5).Value = 45.3
$ws.Cells.Item(16,6).Value = 95540
$ws.Cells.Item(16,7).Value = 5.5
$ws.Cells.Item(16,8).Value = 11614
$ws.Cells.Item(16,9).Value = 0.8
$ws.Cells.Item(16,10).Value = 1681
$ws.Cells.Item(16,11).Value = 51.6
$ws.Cells.Item(16,12).Value = 108835
$ws.Cells.Item(16,13).Value = "'3"

# Row 17: MMR016
$ws.Cells.Item(17,1).Value = "MMR016"
$ws.Cells.Item(17,2).Value = 70760
$ws.Cells.Item(17,3).Value = 85.5
$ws.Cells.Item(17,4).Value = 60490
$ws.Cells.Item(17,5).Value = 10.9
$ws.Cells.Item(17,6).Value = 7743
$ws.Cells.Item(17,7).Value = 3.6
$ws.Cells.Item(17,8).Value = 2528
$ws.Cells.Item(17,9).Value = 0
$ws.Cells.Item(17,10).Value = 0
$ws.Cells.Item(17,11).Value = 14.5
$ws.Cells.Item(17,12).Value = 10271
$ws.Cells.Item(17,13).Value = "1-2"

# Row 18: MMR017
$ws.Cells.Item(18,1).Value = "MMR017"
$ws.Cells.Item(18,2).Value = 81840
$ws.Cells.Item(18,3).Value = 92.5
$ws.Cells.Item(18,4).Value = 75687
$ws.Cells.Item(18,5).Value = 3.6
$ws.Cells.Item(18,6).Value = 2975
$ws.Cells.Item(18,7).Value = 3.3
$ws.Cells.Item(18,8).Value = 2695
$ws.Cells.Item(18,9).Value = 0.6
$ws.Cells.Item(18,10).Value = 484
$ws.Cells.Item(18,11).Value = 7.5
$ws.Cells.Item(18,12).Value = 6153
$ws.Cells.Item(18,13).Value = "1-2"

# Row 19: MMR018
$ws.Cells.Item(19,1).Value = "MMR018"
$ws.Cells.Item(19,2).Value = 85275
$ws.Cells.Item(19,3).Value = 87.09999999999999
$ws.Cells.Item(19,4).Value = 74295
$ws.Cells.Item(19,5).Value = 8.5
$ws.Cells.Item(19,6).Value = 7273
$ws.Cells.Item(19,7).Value = 4.3
$ws.Cells.Item(19,8).Value = 3707
$ws.Cells.Item(19,9).Value = 0
$ws.Cells.Item(19,10).Value = 0
$ws.Cells.Item(19,11).Value = 12.9
$ws.Cells.Item(19,12).Value = 10980
$ws.Cells.Item(19,13).Value = "1-2"

# Drop now-unused rows 20:22 (finer admin3 rows no longer needed post-aggregation)
$ws.Rows("20:22").Delete() | Out-Null

Write-Output "done"
